# Update gh-pages to output generated at 456a3b4
#
# The upstream scraper re-ran and the "南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini"
# event (2024-10-03) dropped out of the list entirely, so its row is removed
# and every row below it shifts up by one. Several other rows also picked up
# refreshed vote counts / prices (columns F / G) from the re-scrape, and one
# ticket (the 2024良牙动漫秋季盛典 show) became unavailable for sale.
#
# Both the "展览" sheet and the "全部类型" sheet carry the same event list
# (the latter simply has two extra rows for non-展览 categories), so the same
# edit is applied to each.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the row 6 event ("南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini") —
    # everything below shifts up automatically (dimension shrinks too).
    $ws.Rows.Item(6).Delete()

    # Refresh "想去人数" (F) for rows 3-4, and both "想去人数"/"最低票价"
    # (F/G) for row 5, which has gone off-sale.
    $ws.Cells.Item(3, 6).Value = 34
    $ws.Cells.Item(4, 6).Value = 229
    $ws.Cells.Item(5, 6).Value = 3822
    $ws.Cells.Item(5, 7).Value = "不可售"

    # Renumber the index column (A) now that a row was removed.
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # The last two surviving rows (the 万圣派对 and 万圣漫控嘉年华10 events)
    # also picked up refreshed "想去人数"/"最低票价" numbers.
    $ws.Cells.Item($lastRow - 1, 6).Value = 23
    $ws.Cells.Item($lastRow - 1, 7).Value = 60
    $ws.Cells.Item($lastRow, 6).Value = 434
    $ws.Cells.Item($lastRow, 7).Value = 50
}
